# find-artek_feature_registration.xlsx hotfix
# - Extend the SRID list sheet with more UTM zones (18N-29N) and region notes
# - Point the "Features" sheet's SRID column (F) validation at the new SRID list
#   table instead of the old hard-coded "3182,4326,32622" list
# - Move the active selection on the Features sheet from G2 to F2

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "SRID list" sheet - add the extra SRID/zone rows
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("SRID list")

# Remove the now-superfluous blank row so the note block shifts from row 24
# up to row 23 (6 "#" rows are kept below the data instead of 18 blank ones).
$ws2.Rows.Item(22).Delete()

# Column C (Description) for all the new UTM-zone rows, entered top-to-bottom.
$ws2.Range("C4").Value = "WGS 84 / UTM zone 18N"
$ws2.Range("C5").Value = "WGS 84 / UTM zone 19N"
$ws2.Range("C6").Value = "WGS 84 / UTM zone 20N"
$ws2.Range("C7").Value = "WGS 84 / UTM zone 21N"
$ws2.Range("C8").Value = "WGS 84 / UTM zone 22N"
$ws2.Range("C9").Value = "WGS 84 / UTM zone 23N"
$ws2.Range("C10").Value = "WGS 84 / UTM zone 24N"
$ws2.Range("C11").Value = "WGS 84 / UTM zone 25N"
$ws2.Range("C12").Value = "WGS 84 / UTM zone 26N"
$ws2.Range("C13").Value = "WGS 84 / UTM zone 27N"
$ws2.Range("C14").Value = "WGS 84 / UTM zone 28N"
$ws2.Range("C15").Value = "WGS 84 / UTM zone 29N"

# Column B (SRID code) for the same rows.
$ws2.Range("B4").Value = 32618
$ws2.Range("B5").Value = 32619
$ws2.Range("B6").Value = 32620
$ws2.Range("B7").Value = 32621
$ws2.Range("B8").Value = 32622
$ws2.Range("B9").Value = 32623
$ws2.Range("B10").Value = 32624
$ws2.Range("B11").Value = 32625
$ws2.Range("B12").Value = 32626
$ws2.Range("B13").Value = 32627
$ws2.Range("B14").Value = 32628
$ws2.Range("B15").Value = 32629

# Column A (#) sequence number for the data rows.
$ws2.Range("A4").Value = 3
$ws2.Range("A5").Value = 4
$ws2.Range("A6").Value = 5
$ws2.Range("A7").Value = 6
$ws2.Range("A8").Value = 7
$ws2.Range("A9").Value = 8
$ws2.Range("A10").Value = 9
$ws2.Range("A11").Value = 10
$ws2.Range("A12").Value = 11
$ws2.Range("A13").Value = 12
$ws2.Range("A14").Value = 13
$ws2.Range("A15").Value = 14

# Column D (Comment) - region notes. Entered with the repeated "East Greenland"
# value first, then the remaining unique call-outs.
$ws2.Range("D9").Value = "South Greenland (Narsarssuaq, Nanortalik, etc.)"
$ws2.Range("D10").Value = "East Greenland"
$ws2.Range("D11").Value = "East Greenland"
$ws2.Range("D12").Value = "East Greenland"
$ws2.Range("D13").Value = "East Greenland"
$ws2.Range("D14").Value = "East Greenland"
$ws2.Range("D15").Value = "East Greenland"
$ws2.Range("D7").Value = "Upernavik area"
$ws2.Range("D5").Value = "Thule / Qaanaaq area"
$ws2.Range("D8").Value = "Most of West Greenland covered"
$ws2.Range("D4").Value = ""
$ws2.Range("D6").Value = ""

# Continue the "#" numbering below the data table (rows 16-21), leaving the
# remaining columns blank.
$ws2.Range("A16").Value = 15
$ws2.Range("A17").Value = 16
$ws2.Range("A18").Value = 17
$ws2.Range("A19").Value = 18
$ws2.Range("A20").Value = 19
$ws2.Range("A21").Value = 20

# ---------------------------------------------------------------------------
# 2. "Features" sheet - point the SRID (F) validation at the SRID list table
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Features")

$fRange = $ws1.Range("F2:F1048576")
$fRange.Validation.Delete()
$fRange.Validation.Add(3, 1, 1, "='SRID list'!`$B`$2:`$B`$15")
$fRange.Validation.IgnoreBlank = $true
$fRange.Validation.InCellDropdown = $true
$fRange.Validation.ShowInput = $true
$fRange.Validation.ShowError = $false

# Move the active selection from G2 to F2.
$ws1.Range("F2").Select()

Write-Output "SRID list extended and Features SRID validation repointed."
